## Added PER material and type
## Adds a new "Tube PER" worksheet (after "Tube Acier") containing the PER
## pipe dimensions/cross-section data, makes it the active sheet, and moves
## the selection on "Tube Acier" to C1 (matching the upstream commit).

$wb = $excel.ActiveWorkbook

# --- Add the new "Tube PER" worksheet as the last tab ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Tube PER"

# --- Populate data: inner diameter (A), outer diameter (B), cross-section
#     area formula (C) = (A^2 * PI()) / 4 --------------------------------
$data = @(
    @(9.8, 12),
    @(13, 16),
    @(16.2, 20),
    @(20.4, 25),
    @(26.2, 32)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $data[$i][0]
    $newSheet.Cells.Item($row, 2).Value = $data[$i][1]
    $newSheet.Cells.Item($row, 3).Formula = "=(A$row*A$row*PI())/4"
}

# Page setup matching the other data sheets (A4, portrait)
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- Move the selection on "Tube Acier" to C1 ---------------------------
$ws3 = $wb.Worksheets.Item("Tube Acier")
$ws3.Activate()
$ws3.Range("C1").Select()

# --- Finish with "Tube PER" active/selected (D3), as the last edited tab
$newSheet.Activate()
$newSheet.Range("D3").Select()
